{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Goal (per the diff): append a new paragraph at the very end of the\n// document body (right before the final section break), reading:\n//   \"4. To do that the buttons in the message window will appear only for commite member.\"\n// with \"commite\" flagged by the spell-checker proof marks (<w:proofErr .../>),\n// and the paragraph carrying <w:bidi w:val=\"0\"/> like its neighboring paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The new paragraph goes after the current last paragraph in the body\n// (i.e. right before the closing sectPr), matching the diff which inserts\n// it right after the \"...more correct.\" paragraph.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst lastRange = lastParagraph.getRange();\n\n// Flat OPC (pkg:package) fragment describing the single new paragraph,\n// including the spell-check proof marks around \"commite\" exactly as the\n// diff specifies, and the <w:bidi w:val=\"0\"/> paragraph property that\n// matches the rest of the document's paragraphs.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:bidi w:val=\"0\"/>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">4. To do that the buttons in the message window will appear only for </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>commite</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> member.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nlastRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $doc / $app resolve to the running Application/Document; the\n# target document is $word.ActiveDocument.\n#\n# Goal (per the diff): append a new paragraph at the very end of the\n# document body (right before the final section break), reading:\n#   \"4. To do that the buttons in the message window will appear only for commite member.\"\n# with \"commite\" flagged by the spell-checker proof marks (<w:proofErr .../>),\n# and the paragraph carrying <w:bidi w:val=\"0\"/> like its neighboring paragraphs.\n\n$d = $word.ActiveDocument\n\n# Locate the current last paragraph (\"...more correct.\") so the new\n# paragraph lands right after it and before the closing section break.\n$lastParagraph = $d.Paragraphs.Last\n$lastRange = $lastParagraph.Range\n$lastRange.InsertParagraphAfter()\n\n# The freshly created (empty) paragraph is now the document's last one.\n$newParagraph = $d.Paragraphs.Last\n$newRange = $newParagraph.Range\n\n# Flat WordprocessingML fragment for the whole paragraph, including the\n# spell-check proof marks around \"commite\" exactly as the diff specifies,\n# and the <w:bidi w:val=\"0\"/> paragraph property matching the rest of the\n# document's paragraphs. InsertXML replaces the (empty) target range's\n# paragraph contents with this markup.\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n       '<w:pPr><w:bidi w:val=\"0\"/></w:pPr>' +\n       '<w:r><w:t xml:space=\"preserve\">4. To do that the buttons in the message window will appear only for </w:t></w:r>' +\n       '<w:proofErr w:type=\"spellStart\"/>' +\n       '<w:r><w:t>commite</w:t></w:r>' +\n       '<w:proofErr w:type=\"spellEnd\"/>' +\n       '<w:r><w:t xml:space=\"preserve\"> member.</w:t></w:r>' +\n       '</w:p>'\n\n$newRange.InsertXML($xml)\n"}
